$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Title: "Check-Out Patient" -> "Delete Database", and move the
#    "_GoBack" bookmark so it sits right after the new title run
#    (this mirrors Word's own "last-edit" bookmark housekeeping).
#    We rebuild the whole paragraph (text + bookmark) via InsertXML so
#    the bookmark lands inline, in the same paragraph, right after the
#    run - not as a brand new paragraph.
# ---------------------------------------------------------------------
$titlePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*Check-Out Patient*") {
        $titlePara = $p
        break
    }
}

if ($titlePara -ne $null) {
    $titleXml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main' w:rsidR='00794829' w:rsidRDefault='0029069A' w:rsidP='00794829'><w:pPr><w:widowControl w:val='0'/><w:autoSpaceDE w:val='0'/><w:autoSpaceDN w:val='0'/><w:adjustRightInd w:val='0'/><w:jc w:val='center'/><w:rPr><w:rFonts w:ascii='Verdana' w:hAnsi='Verdana'/><w:b/><w:bCs/><w:sz w:val='32'/><w:szCs w:val='52'/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii='Verdana' w:hAnsi='Verdana'/><w:b/><w:bCs/><w:sz w:val='32'/><w:szCs w:val='52'/></w:rPr><w:t>Delete Database</w:t></w:r><w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/></w:p>"
    $titlePara.Range.InsertXML($titleXml) | Out-Null
}

# ---------------------------------------------------------------------
# 2) Merge the two runs describing pre-conditions ("The database needs
#    to be updated ... logged in" + ". The board member has already
#    selected mentor/mentees. ") into a single run/text node. A no-op
#    Find/Replace over the full combined text coalesces the two runs.
# ---------------------------------------------------------------------
$preReqText = "The database needs to be updated to the latest version and the board member needs to be logged in. The board member has already selected mentor/mentees. "
$d.Content.Find.Execute($preReqText, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $preReqText, 2) | Out-Null

# ---------------------------------------------------------------------
# 3) Remove the stray "_GoBack" bookmark that used to sit after
#    "...download database" button" (it has effectively moved to the
#    title paragraph above).
# ---------------------------------------------------------------------
$dbPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*download database*button*") {
        $dbPara = $p
        break
    }
}

if ($dbPara -ne $null) {
    $dbXml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main' w:rsidR='0016302B' w:rsidRPr='0029069A' w:rsidRDefault='000A1BD4' w:rsidP='0029069A'><w:pPr><w:widowControl w:val='0'/><w:autoSpaceDE w:val='0'/><w:autoSpaceDN w:val='0'/><w:adjustRightInd w:val='0'/><w:rPr><w:rFonts w:ascii='Verdana' w:hAnsi='Verdana'/><w:sz w:val='18'/><w:szCs w:val='18'/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/><w:sz w:val='20'/></w:rPr><w:t>The board members clicks the &#8220;download database&#8221; button</w:t></w:r></w:p>"
    $dbPara.Range.InsertXML($dbXml) | Out-Null
}

# ---------------------------------------------------------------------
# 4) Merge the two runs "The " + "database CSV is not corrupted and
#    can be opened in the board member's computer" into a single run.
# ---------------------------------------------------------------------
$csvText = "The database CSV is not corrupted and can be opened in the board member" + [char]0x2019 + "s computer"
$d.Content.Find.Execute($csvText, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $csvText, 2) | Out-Null
